$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "José"
$ws.Range("B2").Value = "300.000.003-12"
$ws.Range("C2").Value = "ricardo foda"
$ws.Range("D2:K2").ClearContents()

$ws.Range("D2").Select()
